$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.905.36'
$ws.Range('E2').Value = '  +0.68%  '

$ws.Range('D3').Value = '2.033.59'
$ws.Range('E3').Value = '  -0.29%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '227.93'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.90%  '

$ws.Range('E6').Value = '  -0.15%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '60.48'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +7.25%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.385'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.36%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0811'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.71%  '

$ws.Range('E11').Value = '  +0.96%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.67'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.42%  '

$ws.Range('D13').Value = '2.335.21'
$ws.Range('E13').Value = '  -0.19%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '21.28'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +4.18%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.757'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.73%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.25'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.32%  '

$ws.Range('D17').Value = '2.040.78'
$ws.Range('E17').Value = '  -0.07%  '

$ws.Range('D18').Value = '37.929.68'
$ws.Range('E18').Value = '  +0.95%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.07'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.24%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.68'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.80%  '

$ws.Range('D21').Value = '0.0₃0825'
$ws.Range('E21').Value = '  -0.25%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '224.78'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.58%  '

$ws.Range('E23').Value = '  +0.09%  '

$ws.Range('E24').Value = '  -1.09%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.21'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.37%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '165.53'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.10%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.19'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.24%  '

$ws.Range('E28').Value = '  -2.43%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.90'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.67%  '

$ws.Range('E30').Value = '  -3.39%  '

$ws.Range('E31').Value = '  +1.62%  '

$ws.Range('E32').Value = '  -0.99%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.05'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.34%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.52'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.56%  '

$ws.Range('E35').Value = '  -0.85%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.35'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +9.26%  '

$ws.Range('E37').Value = '  -3.17%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.25'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.13%  '

$ws.Range('E39').Value = '  -0.03%  '

$ws.Range('D40').Value = '1.529.95'
$ws.Range('E40').Value = '  +3.43%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0218'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.22%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '96.93'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.66%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '16.59'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.27%  '

$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0922'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.93%  '

$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.78'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.12%  '

$ws.Range('E46').Value = '  -0.28%  '

$ws.Range('E47').Value = '  -3.34%  '

$ws.Range('E48').Value = '  +1.09%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.01'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.65%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.05'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.15%  '

$ws.Range('D51').Value = '2.223.53'
$ws.Range('E51').Value = '  -0.08%  '
